$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the grades for row 14 (Кузнецов Владимир) columns C:F
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 5

# Update the active selection to match the final cursor position
$ws.Range("P16").Select()
